$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E1").Value = "lat"
$ws.Range("F1").Value = "lng"

$ws.Range("E2").Value = 51.5074
$ws.Range("F2").Value = 0.1278

$ws.Range("E3").Formula = "=E2+0.001"
$ws.Range("F3").Formula = "=F2-0.001"

$ws.Range("E4:E12").Formula = "=E3+0.001"
$ws.Range("F4:F12").Formula = "=F3-0.001"

$ws.Range("G3").Select()

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
